$wb = $excel.ActiveWorkbook

# --- "Range Status" sheet: B2:B7 reset to 0, C2:C7 cleared entirely ---
$wsRange = $wb.Worksheets.Item("Range Status")
$wsRange.Range("B2:B7").Value = 0
$wsRange.Range("C2:C7").ClearContents()

# --- "Species qualification" sheet: B5 (Range Analysis) reset to 0 ---
$wsQual = $wb.Worksheets.Item("Species qualification")
$wsQual.Range("B5").Value = 0
